$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): the "Factor Quartile" spans are narrowing from 5
# columns per market (quintiles 0-4) to 4 columns per market (quartiles 0-3),
# so the "Emerging" label moves from G1 to F1 and the merged ranges shrink.

# Unmerge the existing header ranges so individual cells can be edited
$ws.Range("B1:F1").UnMerge()
$ws.Range("G1:K1").UnMerge()

# Move the "Emerging" label from G1 to F1
$ws.Range("F1").Value = "Emerging"
$ws.Range("G1").ClearContents()

# Re-merge the header cells over the new (narrower) column ranges
$ws.Range("B1:E1").Merge()
$ws.Range("F1:I1").Merge()

# Merging re-splits the per-cell border formatting; restore the uniform
# bordered look so every row-1 cell keeps using the original style.
$ws.Range("A1:K1").Borders.LineStyle = 1
$ws.Range("A1:K1").Borders.ColorIndex = -4105

# --- Row 2 (quartile index labels): shift values left into the narrower
# column layout, now only 4 buckets (0-3) per market instead of 5 (0-4).
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 3

# Remove the now-unused columns J and K entirely (row 1 and row 2)
$ws.Range("J1:K1").Clear()
$ws.Range("J2:K2").Clear()

# --- Row 4: newly (re)computed forward-return values
$ws.Range("B4").Value = 0.008944568700538868
$ws.Range("C4").Value = 0.005366422004796994
$ws.Range("D4").Value = 0.006905512689471749
$ws.Range("E4").Value = 0.007906866970260331
$ws.Range("F4").Value = 0.0202342233071515
$ws.Range("G4").Value = 0.01638605994695274
$ws.Range("H4").Value = 0.007163406606804336
$ws.Range("I4").Value = 0.005088274823469566

# Remove the now-unused columns J and K entirely (row 4)
$ws.Range("J4:K4").Clear()
